# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Rewrites the worker/period detail rows (16-38) of the account-statement
# sheet, grouping the rows by worker (instead of by period) and updating
# the overdue values accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, DocType, DocNumber, WorkerName, Period, OverdueValue
$rows = @(
    @(16, "CC", "10967447",   "FRANCISCO AMADOR GUERRA GUERRERO", "1812", 28124),
    @(17, "CC", "10967447",   "FRANCISCO AMADOR GUERRA GUERRERO", "1902", 31249),
    @(18, "CC", "10967447",   "FRANCISCO AMADOR GUERRA GUERRERO", "1811", 31249),
    @(19, "CC", "10967447",   "FRANCISCO AMADOR GUERRA GUERRERO", "1810", 31249),
    @(20, "CC", "10967447",   "FRANCISCO AMADOR GUERRA GUERRERO", "1809", 31249),
    @(21, "CC", "10967447",   "FRANCISCO AMADOR GUERRA GUERRERO", "1807", 31249),
    @(22, "CC", "10967447",   "FRANCISCO AMADOR GUERRA GUERRERO", "1806", 31249),
    @(23, "CC", "1049564184", "WILBERTO TORRES PEREZ",            "1902", 28124),
    @(24, "CC", "1049564184", "WILBERTO TORRES PEREZ",            "1812", 31249),
    @(25, "CC", "1049564184", "WILBERTO TORRES PEREZ",            "1811", 31249),
    @(26, "CC", "1049564184", "WILBERTO TORRES PEREZ",            "1810", 31249),
    @(27, "CC", "1049564184", "WILBERTO TORRES PEREZ",            "1809", 31249),
    @(28, "CC", "1049564184", "WILBERTO TORRES PEREZ",            "1807", 31249),
    @(29, "CC", "1049564184", "WILBERTO TORRES PEREZ",            "1806", 31249),
    @(30, "CC", "7922859",    "JESUS ANTONIO MENDOZA GUERRERO",   "1807", 31249),
    @(31, "CC", "7922859",    "JESUS ANTONIO MENDOZA GUERRERO",   "1806", 31249),
    @(32, "CC", "91077212",   "SAMUEL LEON SUAREZ",               "1902", 28124),
    @(33, "CC", "91077212",   "SAMUEL LEON SUAREZ",               "1812", 31249),
    @(34, "CC", "91077212",   "SAMUEL LEON SUAREZ",               "1811", 31249),
    @(35, "CC", "91077212",   "SAMUEL LEON SUAREZ",               "1810", 31249),
    @(36, "CC", "91077212",   "SAMUEL LEON SUAREZ",               "1809", 31249),
    @(37, "CC", "91077212",   "SAMUEL LEON SUAREZ",               "1807", 31249),
    @(38, "CC", "91077212",   "SAMUEL LEON SUAREZ",               "1806", 5208)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("B$rowNum").Value = $r[1]
    $ws.Range("C$rowNum").Value = $r[2]
    $ws.Range("D$rowNum").Value = $r[3]
    $ws.Range("E$rowNum").Value = $r[4]
    $ws.Range("F$rowNum").Value = $r[5]
}
